$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A11"  = -21.086
    "B11"  = 7.061
    "A12"  = -21.526
    "A15"  = -21.819
    "B23"  = 7.507000000000001
    "A27"  = -21.522
    "A28"  = -20.955
    "B28"  = 6.526000000000001
    "A31"  = -21.849
    "A32"  = -21.293
    "B32"  = 7.375999999999999
    "B34"  = 6.74
    "A36"  = -21.318
    "B36"  = 6.885
    "B37"  = 7.401999999999999
    "A38"  = -20.443
    "B42"  = 8.046000000000001
    "A46"  = -21.615
    "B49"  = 6.691
    "A54"  = -21.108
    "B54"  = 5.524
    "A55"  = -22.184
    "A56"  = -21.461
    "A67"  = -21.422
    "A69"  = -21.422
    "A72"  = -21.689
    "A73"  = -20.078
    "B78"  = 7.131
    "B80"  = 7.772999999999999
    "A83"  = -21.065
    "A86"  = -21.859
    "A91"  = -20.79
    "A93"  = -21.435
    "B97"  = 5.191
    "A99"  = -21.252
    "B99"  = 5.483
    "B100" = 5.172
    "B101" = 5.217000000000001
    "A104" = -21.13
    "A105" = -20.756
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
